$d = $word.ActiveDocument

# --- 1) Collapse the split "(ignore those which are already incorporated
#        in your book's version and date)" runs into a single run.
#        Occurs 3 times in the document; wdReplaceAll (2) on $d.Content
#        handles every occurrence in one call.
$rsquo = [char]0x2019
$ignorePhrase = "(ignore those which are already incorporated in your book" + $rsquo + "s version and date)"
$d.Content.Find.Execute($ignorePhrase, $true, $false, $false, $false, $false, $true, 1, $false, $ignorePhrase, 2) | Out-Null

# --- 2) Collapse the split "(no elision for "a"" runs into a single run.
$ldq = [char]0x201C
$rdq = [char]0x201D
$elisionPhrase = "(no elision for " + $ldq + "a" + $rdq
$d.Content.Find.Execute($elisionPhrase, $true, $false, $false, $false, $false, $true, 1, $false, $elisionPhrase, 2) | Out-Null

# --- 3) Restructure the first summary table: merge the "Section, Paragraph /
#        Reference" column into the "As Printed" column (drop the old first
#        column, keep the second column's content/formatting which already
#        reads "As Printed" / holds the Latha-font empty paragraph mark),
#        then grow the two remaining columns to fill the freed width.
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Delete()

$t.PreferredWidthType = 3
$t.PreferredWidth = 14679 / 20

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $c1 = $t.Cell($r, 1)
    $c1.PreferredWidthType = 3
    $c1.PreferredWidth = 6883 / 20
    $c1.Width = 6883 / 20

    $c2 = $t.Cell($r, 2)
    $c2.PreferredWidthType = 3
    $c2.PreferredWidth = 7796 / 20
    $c2.Width = 7796 / 20
}
